$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.892.03'
$ws.Range("E2").Value = '  +2.90%  '
$ws.Range("D3").Value = '2.611.78'
$ws.Range("E3").Value = '  +1.37%  '
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("E5").Value = '  +0.48%  '
$ws.Range("D6").Value = '''143.57'
$ws.Range("E6").Value = '  +0.22%  '
$ws.Range("E7").Value = '  -0.23%  '
$ws.Range("E8").Value = '  +1.07%  '
$ws.Range("D9").Value = '2.637.63'
$ws.Range("E9").Value = '  +2.16%  '
$ws.Range("E10").Value = '  -2.35%  '
$ws.Range("E11").Value = '  +3.08%  '
$ws.Range("E12").Value = '  -3.13%  '
$ws.Range("E13").Value = '  +6.86%  '
$ws.Range("D14").Value = '3.101.15'
$ws.Range("E14").Value = '  +2.54%  '
$ws.Range("D15").Value = '60.915.51'
$ws.Range("E15").Value = '  +2.91%  '
$ws.Range("D16").Value = '''23.55'
$ws.Range("E16").Value = '  +4.88%  '
$ws.Range("E17").Value = '  +3.18%  '
$ws.Range("D18").Value = '2.622.21'
$ws.Range("E18").Value = '  +1.51%  '
$ws.Range("B19").Value = 'Polkadot'
$ws.Range("C19").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D19").Value = '''4.68'
$ws.Range("E19").Value = '  +3.38%  '
$ws.Range("B20").Value = 'Chainlink'
$ws.Range("C20").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D20").Value = '''11.26'
$ws.Range("E20").Value = '  +9.90%  '
$ws.Range("D21").Value = '''350.25'
$ws.Range("E21").Value = '  +3.46%  '
$ws.Range("D22").Value = '''7.16'
$ws.Range("E22").Value = '  +14.68%  '
$ws.Range("E23").Value = '  +0.37%  '
$ws.Range("D24").Value = '''0.523'
$ws.Range("E24").Value = '  +14.10%  '
$ws.Range("D25").Value = '''64.13'
$ws.Range("E25").Value = '  -0.58%  '
$ws.Range("D26").Value = '''0.163'
$ws.Range("E26").Value = '  +1.48%  '
$ws.Range("D27").Value = '''0.997'
$ws.Range("E27").Value = '  +0.40%  '
$ws.Range("D28").Value = '''7.71'
$ws.Range("E28").Value = '  +6.29%  '
$ws.Range("E29").Value = '  +1.91%  '
$ws.Range("E30").Value = '  +7.28%  '
$ws.Range("D31").Value = '''0.998'
$ws.Range("E31").Value = '  -0.02%  '
$ws.Range("E32").Value = '  +4.25%  '
$ws.Range("D33").Value = '''160.31'
$ws.Range("E33").Value = '  +1.00%  '
$ws.Range("D34").Value = '''19.51'
$ws.Range("E34").Value = '  +2.62%  '
$ws.Range("D35").Value = '''4.27'
$ws.Range("D36").Value = '''0.960'
$ws.Range("E36").Value = '  +10.00%  '
$ws.Range("E37").Value = '  +4.73%  '
$ws.Range("E38").Value = '  +6.24%  '
$ws.Range("D39").Value = '''37.77'
$ws.Range("E39").Value = '  +1.58%  '
$ws.Range("D40").Value = '''0.856'
$ws.Range("E40").Value = '  -1.82%  '
$ws.Range("D42").Value = '''299.15'
$ws.Range("E42").Value = '  +1.86%  '
$ws.Range("D43").Value = '''140.31'
$ws.Range("E43").Value = '  +9.42%  '
$ws.Range("D44").Value = '''0.0989'
$ws.Range("E44").Value = '  +1.27%  '
$ws.Range("E45").Value = '  -0.33%  '
$ws.Range("D46").Value = '''0.607'
$ws.Range("E46").Value = '  +2.27%  '
$ws.Range("D47").Value = '''0.0548'
$ws.Range("E47").Value = '  +2.25%  '
$ws.Range("D48").Value = '''0.0242'
$ws.Range("E48").Value = '  +3.81%  '
$ws.Range("B49").Value = 'WhiteBITCoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D49").Value = '''10.69'
$ws.Range("E49").Value = '  +0.51%  '
$ws.Range("B50").Value = 'InjectiveProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D50").Value = '''19.72'
$ws.Range("E50").Value = '  +6.75%  '
$ws.Range("E51").Value = '  +7.86%  '
